{"js": "// Replace each two-digit multiplication equation with its updated answer.\n// All strings below are unique within the document, so a direct\n// search-and-replace keyed on the old text is safe and order independent.\nconst replacements = [\n  [\"23\u00d779=1817\", \"51\u00d739=1989\"],\n  [\"15\u00d783=1245\", \"36\u00d747=1692\"],\n  [\"78\u00d774=5772\", \"82\u00d731=2542\"],\n  [\"57\u00d713=741\", \"98\u00d740=3920\"],\n  [\"93\u00d763=5859\", \"60\u00d786=5160\"],\n  [\"80\u00d724=1920\", \"41\u00d746=1886\"],\n  [\"76\u00d776=5776\", \"26\u00d772=1872\"],\n  [\"86\u00d792=7912\", \"13\u00d765=845\"],\n  [\"34\u00d716=544\", \"41\u00d725=1025\"],\n  [\"78\u00d733=2574\", \"14\u00d716=224\"],\n  [\"15\u00d791=1365\", \"27\u00d730=810\"],\n  [\"96\u00d716=1536\", \"96\u00d724=2304\"],\n  [\"50\u00d779=3950\", \"96\u00d791=8736\"],\n  [\"80\u00d745=3600\", \"46\u00d768=3128\"],\n  [\"31\u00d714=434\", \"98\u00d792=9016\"],\n  [\"15\u00d715=225\", \"99\u00d726=2574\"],\n  [\"21\u00d713=273\", \"73\u00d732=2336\"],\n  [\"41\u00d733=1353\", \"40\u00d789=3560\"],\n  [\"15\u00d782=1230\", \"19\u00d767=1273\"],\n  [\"32\u00d789=2848\", \"14\u00d781=1134\"],\n  [\"52\u00d777=4004\", \"47\u00d796=4512\"],\n  [\"35\u00d740=1400\", \"57\u00d763=3591\"],\n  [\"28\u00d789=2492\", \"39\u00d768=2652\"],\n  [\"15\u00d788=1320\", \"41\u00d774=3034\"],\n  [\"60\u00d746=2760\", \"24\u00d755=1320\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication equation with its updated answer.\n# All strings below are unique within the document, so a direct\n# Find/Replace (one occurrence each) keyed on the old text is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"23\u00d779=1817\", \"51\u00d739=1989\"),\n    @(\"15\u00d783=1245\", \"36\u00d747=1692\"),\n    @(\"78\u00d774=5772\", \"82\u00d731=2542\"),\n    @(\"57\u00d713=741\", \"98\u00d740=3920\"),\n    @(\"93\u00d763=5859\", \"60\u00d786=5160\"),\n    @(\"80\u00d724=1920\", \"41\u00d746=1886\"),\n    @(\"76\u00d776=5776\", \"26\u00d772=1872\"),\n    @(\"86\u00d792=7912\", \"13\u00d765=845\"),\n    @(\"34\u00d716=544\", \"41\u00d725=1025\"),\n    @(\"78\u00d733=2574\", \"14\u00d716=224\"),\n    @(\"15\u00d791=1365\", \"27\u00d730=810\"),\n    @(\"96\u00d716=1536\", \"96\u00d724=2304\"),\n    @(\"50\u00d779=3950\", \"96\u00d791=8736\"),\n    @(\"80\u00d745=3600\", \"46\u00d768=3128\"),\n    @(\"31\u00d714=434\", \"98\u00d792=9016\"),\n    @(\"15\u00d715=225\", \"99\u00d726=2574\"),\n    @(\"21\u00d713=273\", \"73\u00d732=2336\"),\n    @(\"41\u00d733=1353\", \"40\u00d789=3560\"),\n    @(\"15\u00d782=1230\", \"19\u00d767=1273\"),\n    @(\"32\u00d789=2848\", \"14\u00d781=1134\"),\n    @(\"52\u00d777=4004\", \"47\u00d796=4512\"),\n    @(\"35\u00d740=1400\", \"57\u00d763=3591\"),\n    @(\"28\u00d789=2492\", \"39\u00d768=2652\"),\n    @(\"15\u00d788=1320\", \"41\u00d774=3034\"),\n    @(\"60\u00d746=2760\", \"24\u00d755=1320\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n}\n"}
